# Commit: "add pdf for sl slides"
#
# For slide 1, the speaker notes text box ("Notes Placeholder 2") had all
# of its paragraphs removed:
#   - "Explain R6 in terms of how they will interact with it, eg
#      Lrnr_rf$train vs randomForest"
#   - "what you can do with trained learners, eg predict, importance, "
#   - "S4 methods bundled at class level, eg predict.randomForest,
#      instead of class level"
# leaving the notes text box empty.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$notesPage = $s.NotesPage
for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
    $shp = $notesPage.Shapes.Item($i)
    if ($shp.Name -eq "Notes Placeholder 2" -and $shp.HasTextFrame) {
        $shp.TextFrame.TextRange.Text = ""
    }
}
